# Generate Report for Handoff
# Updates the per-file localization status rows (041e2fed, 404c7296,
# 7d99c921, ed094cd0 -> rows 4-7) on the zh-cn and de-de sheets:
#   - Priority goes from "low" to "ht" (the handoff work was generated)
#   - Latest Handoff Datetime is refreshed to the new generation time

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7, column E = Priority, column H = Latest Handoff Datetime
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-19 04:27:27"
}

# de-de sheet: rows 4-7, column E = Priority, column H = Latest Handoff Datetime
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-19 04:27:32"
}

# Overview sheet: rows 4-7, column G = Latest HO Xliff Generate Date
# (mirrors the de-de handoff timestamp refresh above)
foreach ($r in 4..7) {
    $overview.Cells.Item($r, 7).Value = "2016-08-19 04:27:32"
}
